# Update Daily Report: 2026-01-16
# Appends the 2026-01-15 (serial 46037) depository rows to Daily_Data,
# and refreshes the dependent Today_Summary / Monthly_Stats rollups.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append 22 new rows (200-221) for date 46037 (2026-01-15)
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Daily_Data")

$newData = New-Object 'object[,]' 22,8
$newData[0,0] = 46037
$newData[0,1] = "ASAHI DEPOSITORY LLC Registered"
$newData[0,2] = 0
$newData[0,3] = 0
$newData[0,4] = 0
$newData[0,5] = 0
$newData[0,6] = 0
$newData[0,7] = 0

$newData[1,0] = 46037
$newData[1,1] = "ASAHI DEPOSITORY LLC Eligible"
$newData[1,2] = 0
$newData[1,3] = 0
$newData[1,4] = 0
$newData[1,5] = 0
$newData[1,6] = 0
$newData[1,7] = 0

$newData[2,0] = 46037
$newData[2,1] = "BRINK'S, INC. Registered"
$newData[2,2] = 95517.77499999999
$newData[2,3] = 0
$newData[2,4] = 0
$newData[2,5] = 0
$newData[2,6] = 0
$newData[2,7] = 95517.77499999999

$newData[3,0] = 46037
$newData[3,1] = "BRINK'S, INC. Eligible"
$newData[3,2] = 23710.274
$newData[3,3] = 0
$newData[3,4] = 0
$newData[3,5] = 0
$newData[3,6] = 0
$newData[3,7] = 23710.274

$newData[4,0] = 46037
$newData[4,1] = "CNT DEPOSITORY, INC. Registered"
$newData[4,2] = 1246.06
$newData[4,3] = 0
$newData[4,4] = 0
$newData[4,5] = 0
$newData[4,6] = 0
$newData[4,7] = 1246.06

$newData[5,0] = 46037
$newData[5,1] = "CNT DEPOSITORY, INC. Eligible"
$newData[5,2] = 0
$newData[5,3] = 0
$newData[5,4] = 0
$newData[5,5] = 0
$newData[5,6] = 0
$newData[5,7] = 0

$newData[6,0] = 46037
$newData[6,1] = "DELAWARE DEPOSITORY Registered"
$newData[6,2] = 1633.941
$newData[6,3] = 0
$newData[6,4] = 0
$newData[6,5] = 0
$newData[6,6] = 0
$newData[6,7] = 1633.941

$newData[7,0] = 46037
$newData[7,1] = "DELAWARE DEPOSITORY Eligible"
$newData[7,2] = 18459.584
$newData[7,3] = 0
$newData[7,4] = 0
$newData[7,5] = 0
$newData[7,6] = 0
$newData[7,7] = 18459.584

$newData[8,0] = 46037
$newData[8,1] = "HSBC BANK, USA Registered"
$newData[8,2] = 1295.223
$newData[8,3] = 0
$newData[8,4] = 0
$newData[8,5] = 0
$newData[8,6] = 0
$newData[8,7] = 1295.223

$newData[9,0] = 46037
$newData[9,1] = "HSBC BANK, USA Eligible"
$newData[9,2] = 9281.978999999999
$newData[9,3] = 99.535
$newData[9,4] = 0
$newData[9,5] = 99.535
$newData[9,6] = 0
$newData[9,7] = 9381.513999999999

$newData[10,0] = 46037
$newData[10,1] = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered"
$newData[10,2] = 2395.448
$newData[10,3] = 0
$newData[10,4] = 0
$newData[10,5] = 0
$newData[10,6] = 0
$newData[10,7] = 2395.448

$newData[11,0] = 46037
$newData[11,1] = "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible"
$newData[11,2] = 0
$newData[11,3] = 0
$newData[11,4] = 0
$newData[11,5] = 0
$newData[11,6] = 0
$newData[11,7] = 0

$newData[12,0] = 46037
$newData[12,1] = "JP MORGAN CHASE BANK NA Registered"
$newData[12,2] = 124991.729
$newData[12,3] = 0
$newData[12,4] = 0
$newData[12,5] = 0
$newData[12,6] = 0
$newData[12,7] = 124991.729

$newData[13,0] = 46037
$newData[13,1] = "JP MORGAN CHASE BANK NA Eligible"
$newData[13,2] = 125407.673
$newData[13,3] = 0
$newData[13,4] = 0
$newData[13,5] = 0
$newData[13,6] = 0
$newData[13,7] = 125407.673

$newData[14,0] = 46037
$newData[14,1] = "LOOMIS INTERNATIONAL (US) LLC Registered"
$newData[14,2] = 68084.33
$newData[14,3] = 0
$newData[14,4] = 0
$newData[14,5] = 0
$newData[14,6] = 0
$newData[14,7] = 68084.33

$newData[15,0] = 46037
$newData[15,1] = "LOOMIS INTERNATIONAL (US) LLC Eligible"
$newData[15,2] = 116365.524
$newData[15,3] = 0
$newData[15,4] = 0
$newData[15,5] = 0
$newData[15,6] = 0
$newData[15,7] = 116365.524

$newData[16,0] = 46037
$newData[16,1] = "MALCA-AMIT USA, LLC Registered"
$newData[16,2] = 395.145
$newData[16,3] = 0
$newData[16,4] = 0
$newData[16,5] = 0
$newData[16,6] = 0
$newData[16,7] = 395.145

$newData[17,0] = 46037
$newData[17,1] = "MALCA-AMIT USA, LLC Eligible"
$newData[17,2] = 0
$newData[17,3] = 0
$newData[17,4] = 0
$newData[17,5] = 0
$newData[17,6] = 0
$newData[17,7] = 0

$newData[18,0] = 46037
$newData[18,1] = "MANFRA, TORDELLA & BROOKES, LLC Registered"
$newData[18,2] = 60301.249
$newData[18,3] = 0
$newData[18,4] = 0
$newData[18,5] = 0
$newData[18,6] = 0
$newData[18,7] = 60301.249

$newData[19,0] = 46037
$newData[19,1] = "MANFRA, TORDELLA & BROOKES, LLC Eligible"
$newData[19,2] = 1068.408
$newData[19,3] = 0
$newData[19,4] = 0
$newData[19,5] = 0
$newData[19,6] = 0
$newData[19,7] = 1068.408

$newData[20,0] = 46037
$newData[20,1] = "STONEX PRECIOUS METALS LLC Registered"
$newData[20,2] = 14122.765
$newData[20,3] = 0
$newData[20,4] = 0
$newData[20,5] = 0
$newData[20,6] = 0
$newData[20,7] = 14122.765

$newData[21,0] = 46037
$newData[21,1] = "STONEX PRECIOUS METALS LLC Eligible"
$newData[21,2] = 16.075
$newData[21,3] = 0
$newData[21,4] = 0
$newData[21,5] = 0
$newData[21,6] = 0
$newData[21,7] = 16.075

$targetRange = $wsData.Range("A200:H221")
$targetRange.Value = $newData

# Match the date-serial display format used by the rest of column A
$wsData.Range("A200:A221").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# 2) Today_Summary: HSBC BANK, USA row (row 6) reflects the new Eligible
#    receipt of 99.535 -> Eligible 9281.979 + 99.535 = 9381.514
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Today_Summary")
$wsSummary.Range("B6").Value = 9381.513999999999
$wsSummary.Range("D6").Value = 10676.737

# ---------------------------------------------------------------------
# 3) Monthly_Stats: 2026-01 grand total (row 2) and the HSBC BANK, USA
#    Eligible monthly detail row (row 15)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")
$wsMonthly.Range("B2").Value = 294409.052
$wsMonthly.Range("D2").Value = 664392.7170000001

$wsMonthly.Range("C15").Value = 99.535
$wsMonthly.Range("E15").Value = 9381.513999999999

